{"js": "// Update the date and each two-digit-by-one-digit division answer cell.\nconst pairs = [\n  [\"2025-12-08 Monday\", \"2025-12-09 Tuesday\"],\n  [\"98\u00f76=16, 2\", \"58\u00f74=14, 2\"],\n  [\"29\u00f73=9, 2\", \"94\u00f73=31, 1\"],\n  [\"32\u00f78=4, 0\", \"87\u00f78=10, 7\"],\n  [\"41\u00f76=6, 5\", \"35\u00f73=11, 2\"],\n  [\"30\u00f74=7, 2\", \"96\u00f79=10, 6\"],\n  [\"59\u00f78=7, 3\", \"20\u00f76=3, 2\"],\n  [\"95\u00f76=15, 5\", \"99\u00f76=16, 3\"],\n  [\"26\u00f75=5, 1\", \"70\u00f78=8, 6\"],\n  [\"12\u00f77=1, 5\", \"74\u00f75=14, 4\"],\n  [\"16\u00f75=3, 1\", \"94\u00f76=15, 4\"],\n  [\"14\u00f78=1, 6\", \"23\u00f77=3, 2\"],\n  [\"69\u00f72=34, 1\", \"34\u00f77=4, 6\"],\n  [\"13\u00f72=6, 1\", \"91\u00f73=30, 1\"],\n  [\"13\u00f77=1, 6\", \"71\u00f78=8, 7\"],\n  [\"34\u00f72=17, 0\", \"22\u00f76=3, 4\"],\n  [\"48\u00f75=9, 3\", \"84\u00f79=9, 3\"],\n  [\"80\u00f79=8, 8\", \"38\u00f75=7, 3\"],\n  [\"72\u00f73=24, 0\", \"25\u00f76=4, 1\"],\n  [\"39\u00f73=13, 0\", \"23\u00f75=4, 3\"],\n  [\"31\u00f76=5, 1\", \"74\u00f72=37, 0\"],\n  [\"37\u00f74=9, 1\", \"43\u00f74=10, 3\"],\n  [\"10\u00f73=3, 1\", \"67\u00f74=16, 3\"],\n  [\"19\u00f79=2, 1\", \"80\u00f72=40, 0\"],\n  [\"41\u00f79=4, 5\", \"36\u00f79=4, 0\"],\n  [\"57\u00f77=8, 1\", \"53\u00f73=17, 2\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();", "ps1": "# Update the date and each two-digit-by-one-digit division answer cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-12-08 Monday\", \"2025-12-09 Tuesday\"),\n    @(\"98\u00f76=16, 2\", \"58\u00f74=14, 2\"),\n    @(\"29\u00f73=9, 2\", \"94\u00f73=31, 1\"),\n    @(\"32\u00f78=4, 0\", \"87\u00f78=10, 7\"),\n    @(\"41\u00f76=6, 5\", \"35\u00f73=11, 2\"),\n    @(\"30\u00f74=7, 2\", \"96\u00f79=10, 6\"),\n    @(\"59\u00f78=7, 3\", \"20\u00f76=3, 2\"),\n    @(\"95\u00f76=15, 5\", \"99\u00f76=16, 3\"),\n    @(\"26\u00f75=5, 1\", \"70\u00f78=8, 6\"),\n    @(\"12\u00f77=1, 5\", \"74\u00f75=14, 4\"),\n    @(\"16\u00f75=3, 1\", \"94\u00f76=15, 4\"),\n    @(\"14\u00f78=1, 6\", \"23\u00f77=3, 2\"),\n    @(\"69\u00f72=34, 1\", \"34\u00f77=4, 6\"),\n    @(\"13\u00f72=6, 1\", \"91\u00f73=30, 1\"),\n    @(\"13\u00f77=1, 6\", \"71\u00f78=8, 7\"),\n    @(\"34\u00f72=17, 0\", \"22\u00f76=3, 4\"),\n    @(\"48\u00f75=9, 3\", \"84\u00f79=9, 3\"),\n    @(\"80\u00f79=8, 8\", \"38\u00f75=7, 3\"),\n    @(\"72\u00f73=24, 0\", \"25\u00f76=4, 1\"),\n    @(\"39\u00f73=13, 0\", \"23\u00f75=4, 3\"),\n    @(\"31\u00f76=5, 1\", \"74\u00f72=37, 0\"),\n    @(\"37\u00f74=9, 1\", \"43\u00f74=10, 3\"),\n    @(\"10\u00f73=3, 1\", \"67\u00f74=16, 3\"),\n    @(\"19\u00f79=2, 1\", \"80\u00f72=40, 0\"),\n    @(\"41\u00f79=4, 5\", \"36\u00f79=4, 0\"),\n    @(\"57\u00f77=8, 1\", \"53\u00f73=17, 2\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
